# weight_tracker: "add update 1 day"
# Appends 6 new weigh-in readings (rows 278-283) to the raw_data sheet,
# matching the layout of the existing rows (date+time in A, time-of-day
# fraction in B, weight in C, an AM/PM lookup formula in D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")
$ws.Activate()

# New rows, in the same (reverse-chronological) order they were added to
# the sheet: date/time serial, time-only fraction, weight (kg), AM/PM.
$newRows = @(
    @{ Row = 278; A = 44120.381249999999; B = 0.38125000000000003; C = 70.5 },
    @{ Row = 279; A = 44120.380555555559; B = 0.38055555555555554; C = 70.5 },
    @{ Row = 280; A = 44120.313888888886; B = 0.31388888888888888; C = 70.5 },
    @{ Row = 281; A = 44119.930555555555; B = 0.93055555555555547; C = 71.2 },
    @{ Row = 282; A = 44119.379861111112; B = 0.37986111111111115; C = 70.5 },
    @{ Row = 283; A = 44119.325694444444; B = 0.32569444444444445; C = 70.8 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value = $r.A
    $cellA.NumberFormat = "m/d/yy h:mm"

    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellB.Value = $r.B
    $cellB.NumberFormat = "h:mm"

    $cellC = $ws.Cells.Item($rowNum, 3)
    $cellC.Value = $r.C

    $cellD = $ws.Cells.Item($rowNum, 4)
    $cellD.Formula = "=IF(B$rowNum<TIME(12,0,0), ""AM"", ""PM"")"
}

# Match the saved selection/scroll state recorded in the diff.
$ws.Range("C286").Select() | Out-Null
